$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.798.80"
$ws.Range("E2").Value = "  -1.76%  "

$ws.Range("D3").Value = "3.428.67"
$ws.Range("E3").Value = "  -2.12%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.23%  "

$ws.Range("D7").Value = "3.427.88"
$ws.Range("E7").Value = "  -2.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.87%  "

$ws.Range("E12").Value = "  -3.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.65%  "

$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "3.976.39"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").Value = "3.429.05"
$ws.Range("E17").Value = "  -1.61%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "575.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.99%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "68.838.64"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("E22").Value = "  -3.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "95.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("E25").Value = "  -3.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.93%  "

$ws.Range("E29").Value = "  -4.42%  "

$ws.Range("E30").Value = "  -4.54%  "

$ws.Range("E31").Value = "  -3.65%  "

$ws.Range("E32").Value = "  -3.64%  "

$ws.Range("E33").Value = "  -6.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "585.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0950"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.94%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0463"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.139"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -14.01%  "

$ws.Range("D43").Value = "3.214.21"
$ws.Range("E43").Value = "  -3.66%  "

$ws.Range("D44").Value = "0.0₃0671"
$ws.Range("E44").Value = "  -8.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.80%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.56%  "

$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.292"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.31%  "

$ws.Range("E48").Value = "  -6.70%  "

$ws.Range("E49").Value = "  -2.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("E51").Value = "  +0.00%  "
